$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "location example" cell F4: change the delimiter from
# commas to pipes ("[12,163.5,90]" -> "[12|163.5|90]")
$ws.Range("F4").Value = "[12|163.5|90]"

# New header cells for the extra int[] / TestValue columns (H) and the
# newly introduced column J -- written in the same order the shared
# strings table picked them up in the source edit
$ws.Range("H4").Value = "[[360.1|12|19],[96|1|56],[45|91.5|60]]"
$ws.Range("J3").Value = "int[]"
$ws.Range("J2").Value = "aaaaaa"
$ws.Range("J4").Value = "1,2,3"

# Move the active selection to J5 to match the saved workbook state
$ws.Range("J5").Select()
